$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.872.82'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.633.59'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.37'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.507'
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0634'
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '1.857.65'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.25'
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').Value = '1.626.43'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = '0.0₃0757'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.58'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '25.877.68'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.70'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.40'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.28'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.45'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.126'
$ws.Range('E27').Value = '  +2.67%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.43'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0500'
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.30'
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.902'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').Value = '1.139.05'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.36'
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.799'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('E44').Value = '  -3.22%  '
$ws.Range('D45').Value = '1.767.04'
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('E46').Value = '  +2.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.29'
$ws.Range('E47').Value = '  +1.73%  '
$ws.Range('E48').Value = '  +3.36%  '
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.67'
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.415'
$ws.Range('E51').Value = '  -0.07%  '
